# Update countries & provincias Spain
# Applies the 23-Abril-2020 11:52 data refresh to the "Pais" sheet:
#  - refreshes the statistics for several countries
#  - two pairs of countries swap ranking position (names + figures move
#    to the neighbouring row)
#  - updates the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 11:52"

# --- España (row 5) --------------------------------------------------------
$ws.Range("B5").Value = 213024
$ws.Range("C5").Value = 4635
$ws.Range("D5").Value = 89250
$ws.Range("E5").Value = 101617
$ws.Range("G5").Value = 440
$ws.Range("H5").Value = 22157

# --- Noruega (row 39) --------------------------------------------------
$ws.Range("E39").Value = 7117
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 189

# --- Australia (row 46) -------------------------------------------------
$ws.Range("B46").Value = 6661
$ws.Range("C46").Value = 12
$ws.Range("D46").Value = 5045
$ws.Range("E46").Value = 1541
$ws.Range("F46").Value = 45

# --- Malasia (row 47) ----------------------------------------------------
$ws.Range("B47").Value = 5603
$ws.Range("C47").Value = 71
$ws.Range("D47").Value = 3542
$ws.Range("E47").Value = 1966
$ws.Range("F47").Value = 42
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 95

# --- Banglades / Finlandia swap ranks (rows 51 & 52) ----------------------
$ws.Range("A51").Value = "Finlandia"
$ws.Range("B51").Value = 4284
$ws.Range("C51").Value = 155
$ws.Range("D51").Value = 2000
$ws.Range("E51").Value = 2135
$ws.Range("F51").Value = 63
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 149

$ws.Range("A52").Value = "Banglades"
$ws.Range("B52").Value = 4186
$ws.Range("C52").Value = 414
$ws.Range("D52").Value = 108
$ws.Range("E52").Value = 3951
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 127

# --- Hungria / Kuwait swap ranks (rows 62 & 63) ----------------------------
$ws.Range("A62").Value = "Kuwait"
$ws.Range("B62").Value = 2399
$ws.Range("C62").Value = 151
$ws.Range("D62").Value = 498
$ws.Range("E62").Value = 1887
$ws.Range("F62").Value = 55
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 14

$ws.Range("A63").Value = "Hungria"
$ws.Range("B63").Value = 2284
$ws.Range("C63").Value = 116
$ws.Range("D63").Value = 390
$ws.Range("E63").Value = 1655
$ws.Range("F63").Value = 61
$ws.Range("G63").Value = 14
$ws.Range("H63").Value = 239

# --- Eslovenia (row 77) ---------------------------------------------------
$ws.Range("B77").Value = 1366
$ws.Range("C77").Value = 13
$ws.Range("D77").Value = 211
$ws.Range("E77").Value = 1076
$ws.Range("F77").Value = 23

# --- Cuba / Afganistan swap ranks (rows 80 & 81) ---------------------------
$ws.Range("A80").Value = "Afganistan"
$ws.Range("B80").Value = 1226
$ws.Range("C80").Value = 50
$ws.Range("D80").Value = 177
$ws.Range("E80").Value = 1009
$ws.Range("F80").Value = 7

$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1189
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 341
$ws.Range("E81").Value = 808
$ws.Range("F81").Value = 16

# --- Niger / Albania swap ranks (rows 98 & 99) ------------------------------
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 663
$ws.Range("C98").Value = 29
$ws.Range("D98").Value = 385
$ws.Range("E98").Value = 251
$ws.Range("F98").Value = 4
$ws.Range("H98").Value = 27

$ws.Range("A99").Value = "Niger"
$ws.Range("B99").Value = 662
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 193
$ws.Range("E99").Value = 447
$ws.Range("F99").Value = 0
$ws.Range("H99").Value = 22
